# Actualizacion Datos Personales 4 nov
# Insert a new "rescatable" student record (LOPEZ SANCHEZ CINTHIA) into the
# "Rescatables" sheet at row 5, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Push rows 5..12 down to 6..13 and open up a blank row 5 for the new record.
$ws.Rows("5").Insert()

$ws.Range("A5").Value = 19330051920163
$ws.Range("B5").Value = "LOPEZ"
$ws.Range("C5").Value = "SANCHEZ"
$ws.Range("D5").Value = "CINTHIA"
$ws.Range("E5").Value = "CIENCIA, TECNOLOGÍA, SOCIEDAD Y VALORES"
$ws.Range("F5").Value = "5ALCM"
$ws.Range("G5").Value = 6
